$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three middle date-block columns (F,G,H,I -> collapse into D,E)
# Original layout: A..I = StudentID, Name, Total, 18-03 Status, 18-03 Time, 19-03 Status, 19-03 Time, 20-03 Status, 20-03 Time
# Target layout:   A..E = StudentID, Name, Total, 20-03 Status, 20-03 Time
$ws.Range("D1:G1").EntireColumn.Delete()

# Update header row
$ws.Range("D1").Value = "20-03-2025 Status"
$ws.Range("E1").Value = "20-03-2025 Time"

# Update data rows
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "P"
$ws.Range("E2").Value = "19:44:44"

$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "P"
$ws.Range("E3").Value = "19:45:07"

$ws.Range("C4").Value = 1
$ws.Range("E4").Value = "17:24:11"

$ws.Range("E5").Value = "19:44:31"

$ws.Range("B6").Value = "M. Nawaz"
$ws.Range("D6").Value = "P"
$ws.Range("E6").Value = "19:44:56"
